$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.642.67'
$ws.Range("E2").Value = '  -2.62%  '
$ws.Range("D3").Value = '1.808.24'
$ws.Range("E3").Value = '  -2.12%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.70'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.613'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.09%  '
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '40.04'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -9.52%  '
$ws.Range("E9").Value = '  +2.78%  '
$ws.Range("E10").Value = '  -2.67%  '
$ws.Range("E11").Value = '  -2.26%  '
$ws.Range("D12").Value = '2.068.73'
$ws.Range("E12").Value = '  -2.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.26'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.68%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.818.01'
$ws.Range("E14").Value = '  -1.72%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.664'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.93%  '
$ws.Range("E16").Value = '  -4.47%  '
$ws.Range("D17").Value = '34.489.09'
$ws.Range("E17").Value = '  -2.99%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.28'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.05%  '
$ws.Range("D19").Value = '0.0₃0782'
$ws.Range("E19").Value = '  -3.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '240.07'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.88'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.67'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.30%  '
$ws.Range("E23").Value = '  +0.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.24'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '173.62'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.12%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.74'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.31'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -3.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.123'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.53%  '
$ws.Range("E29").Value = '  -9.57%  '
$ws.Range("E30").Value = '  +0.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.00'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.69%  '
$ws.Range("E32").Value = '  -3.86%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.92'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -5.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.24'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +9.75%  '
$ws.Range("E35").Value = '  -3.30%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.691'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '90.38'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -5.21%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.34'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +5.18%  '
$ws.Range("D39").Value = '1.332.55'
$ws.Range("E39").Value = '  -1.24%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0191'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -3.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.966'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -5.53%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.44'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -5.61%  '
$ws.Range("E43").Value = '  -1.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.23'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -7.51%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.73'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -3.18%  '
$ws.Range("B46").Value = 'Kaspa'
$ws.Range("C46").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0513'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.67%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.12'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.95%  '
$ws.Range("D48").Value = '1.993.44'
$ws.Range("E48").Value = '  -1.52%  '
$ws.Range("E49").Value = '  +0.22%  '
$ws.Range("E50").Value = '  +3.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '97.85'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -4.74%  '
